# edit.ps1 - applies the "rex commit for push" change:
#   1. Merge the two runs "Rex " + "change it 3" into a single run
#      "Rex change it 3" (happens twice in the document).
#   2. Add a new paragraph "Rex change it for push" right after the
#      second "Rex change it 3" paragraph, moving the _GoBack bookmark
#      from the end of that paragraph onto the end of the new one.
#   3. Remove the now-superfluous trailing empty paragraph.

$d = $word.ActiveDocument

# --- Step 1: merge the split "Rex " / "change it 3" runs -------------------
# Running a Find/Replace whose replacement text equals the match text is a
# reliable way to coalesce adjacent runs that share identical formatting
# into a single run, without touching the visible text. wdReplaceAll (2)
# makes this apply to both occurrences in one call.
$d.Content.Find.Execute("Rex change it 3", $true, $false, $false, $false, $false, $true, 1, $false, "Rex change it 3", 2) | Out-Null

# --- Step 2: insert the new paragraph after the second occurrence ----------
# Paragraph 8 is the second "Rex change it 3" paragraph (the one that also
# carries the _GoBack bookmark).
$p8 = $d.Paragraphs.Item(8)

$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)

# Write the new paragraph's text. A trailing "Z" marker character is added
# temporarily -- placing a bookmark in a zero-length range sitting exactly
# on a paragraph mark is unreliable, so instead a bookmark is wrapped around
# a throw-away character and that character is then deleted via the
# bookmark's own Range, which leaves the (now collapsed) bookmark correctly
# anchored right after the real text.
$rng = $d.Range($p9.Range.Start, $p9.Range.End - 1)
$rng.Text = "Rex change it for pushZ"
$p9 = $d.Paragraphs.Item(9)

$markerRng = $d.Range($p9.Range.End - 2, $p9.Range.End - 1)

# --- Step 3: move the _GoBack bookmark onto the new paragraph --------------
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $markerRng)
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Text = ""

# --- Step 4: drop the trailing empty paragraph that used to precede sectPr -
$p9 = $d.Paragraphs.Item(9)
$last = $d.Paragraphs.Last
$delRng = $d.Range($p9.Range.End - 1, $last.Range.End)
$delRng.Delete()
